$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Failed Tests" figures: Firefox 3.6.6 moves into row 3 (264
# failures) and Jurassic drops to row 4 (253 failures), keeping the
# existing descending sort by failure count intact.
$ws.Range("A3").Value = "Firefox 3.6.6"
$ws.Range("B3").Value = 264
$ws.Range("A4").Value = "Jurassic"
$ws.Range("B4").Value = 253

$ws.Range("E3").Select()
